$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Text = $new
}

function Insert-After($anchor, $text) {
    $rng = $d.Content
    $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Collapse(0)
    $rng.InsertAfter($text)
}

# --- Title ---
Replace-Text "Astronomy's Unseen Universe" "Biology: The Study of Life's Symphony"

# --- Author line: "Dr" + "." + " Neil deGrasse Tyson" -> "Amelia Barton" (single run) ---
Replace-Text "Dr. Neil deGrasse Tyson" "Amelia Barton"

# --- Email ---
Replace-Text "neil@astro" "ameliabarton@gmail"
Replace-Text "nyc" "net"

# --- Body paragraph 1 ---
Replace-Text "In the cosmic tapestry of our existence, humanity has been captivated by the allure of celestial bodies since ancient times" "Biology, the captivating realm of life, has enthralled humankind since time immemorial"
Replace-Text " From the earliest stargazers marveling at the night sky to modern-day cosmologists unraveling the mysteries of the cosmos, astronomy remains a discipline that ignites our imagination and deepens our understanding of the universe" " As we journey through the intricacies of living organisms, we discover a universe teeming with diversity, interconnectedness, and harmonious choreography"
Replace-Text " Yet, tantalizingly, beyond the reach of our telescopes, lies an unseen universe, an enigmatic realm that beckon to be explored" " From the minuscule cells that form the building blocks of life to the awe-inspiring ecosystems that shape our world, biology unravels the secrets of existence"
Insert-After "biology unravels the secrets of existence" ". In this exploration, we'll embark on a quest to understand the fundamentals of life, delving into the remarkable adaptations, processes, and relationships that govern the living world"

Replace-Text "This hidden universe, vast and mysterious, holds secrets of unimaginable proportions" "The study of biology not only broadens our knowledge of the natural world but also challenges us to solve real-world problems"
Replace-Text " It is a cosmic puzzle waiting to be unraveled, beckoning scientists and enthusiasts alike to delve into its enigmatic depths" " As we unravel the mysteries of life, we gain insights into human health, disease, and the intricate balance of ecosystems"
Replace-Text " As we continue to unlock the mysteries of the known universe, the unseen realm holds the promise of new discoveries, transformative theories, and fundamental shifts in our current understanding of reality" " From curing diseases to developing sustainable agricultural practices, biology empowers us to make informed decisions and strive for a healthier future"
Insert-After "strive for a healthier future" ". Moreover, the field offers boundless opportunities for innovation and discovery, beckoning us to unravel the enigmas that still cloak the living world"

Replace-Text "The unseen universe may reveal the nature of dark matter and dark energy, entities that govern the expansion of the cosmos" "Biology is an orchestra of interconnected disciplines, harmoniously weaving together diverse fields to unravel the secrets of life"
Replace-Text " It might harbor undiscovered exoplanets, possibly teeming with life, hidden within the vast expanse of habitable zones" " From genetics to ecology and physiology to evolution, each branch of biology contributes its unique melody to the grand symphony of life"
Replace-Text " The study of gravitational waves and cosmic radiation could provide a window into the cataclysmic events that shaped the universe's creation and evolution" " Moreover, biology seamlessly merges with other disciplines, such as chemistry, physics, and mathematics, creating a tapestry of knowledge that deepens our understanding of the universe"
Insert-After "deepens our understanding of the universe" ". Through this interdisciplinary approach, we gain a holistic perspective, recognizing the unity and interconnectedness of all things"

# --- Summary paragraph ---
Replace-Text "In the celestial tapestry of our cosmos, there exists an unseen universe, a realm shrouded in mystery and wonder" "Delving into the intricacies of life, we uncover the remarkable adaptations, processes, and relationships that govern the living world"
Replace-Text " Concealed from our view, this enigmatic domain holds promises of untapped knowledge, transformative theories, and fundamental shifts in understanding" " Biology not only expands our knowledge of the natural world but also empowers us to solve real-world problems, from curing diseases to preserving ecosystems"
Replace-Text " The exploration of the unseen universe through meticulous study and unwavering curiosity propels us towards a deeper comprehension of our cosmic heritage and place within the vast expanse of existence" " As an orchestra of interconnected disciplines, biology harmoniously weaves together diverse fields, offering boundless opportunities for innovation and discovery"
Insert-After "offering boundless opportunities for innovation and discovery" ". Biology's interdisciplinary nature fosters a holistic understanding of the universe, recognizing the unity and interconnectedness of all things. Through the study of biology, we gain a profound appreciation for the symphony of life and our place within it"

# --- New empty paragraph at end of document body ---
$lastPar = $d.Paragraphs($d.Paragraphs.Count)
$endRng = $lastPar.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
